# Fix mojibake: "Â±" (U+00C2 U+00B1, a double-encoded UTF-8 artifact)
# should be the plain "±" (U+00B1) character, in the parenthesized
# mean/std-dev portion of columns B (f1_score_weighted), C (training_time)
# and D (test_time), for data rows 2 through 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c1 = [char]0x00C2
$c2 = [char]0x00B1
$badChar = "$c1$c2"
$goodChar = "$c2"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = $cell.Value2
        if ($value -ne $null -and $value.Contains($badChar)) {
            $cell.Value = $value.Replace($badChar, $goodChar)
        }
    }
}
